# Reset to 1V8/0V8 setting
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 41 (R64,R38,R30): resistor value 110k -> 115k
$ws.Range("C41").Value = "115k"

# Row 42 (U4): regulator part number MIC5504-MS -> MIC5504-1.8YM5-TR
$ws.Range("C42").Value = "MIC5504-1.8YM5-TR"

# Row 49 (R66,R32,R40): resistor value 220k -> 402k
$ws.Range("C49").Value = "402k"

# Match the saved view/selection state
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C43").Select()
